$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''60.941.72'
$ws.Range("E2").Value = '  -3.19%  '

# Row 3
$ws.Range("D3").Value = '''3.001.35'
$ws.Range("E3").Value = '  -5.44%  '

# Row 4
$ws.Range("E4").Value = '  +0.05%  '

# Row 5
$ws.Range("D5").Value = '''564.68'
$ws.Range("E5").Value = '  -4.37%  '

# Row 6
$ws.Range("D6").Value = '''128.44'
$ws.Range("E6").Value = '  -5.92%  '

# Row 7
$ws.Range("E7").Value = '  +0.06%  '

# Row 8
$ws.Range("D8").Value = '''2.998.85'
$ws.Range("E8").Value = '  -5.42%  '

# Row 9
$ws.Range("E9").Value = '  -2.94%  '

# Row 10
$ws.Range("E10").Value = '  -5.86%  '

# Row 11
$ws.Range("D11").Value = '''5.23'
$ws.Range("E11").Value = '  -1.25%  '

# Row 12
$ws.Range("D12").Value = '''0.433'
$ws.Range("E12").Value = '  -5.29%  '

# Row 13
$ws.Range("D13").Value = '''0.0000223'
$ws.Range("E13").Value = '  -6.49%  '

# Row 14
$ws.Range("D14").Value = '''32.94'
$ws.Range("E14").Value = '  -5.62%  '

# Row 15
$ws.Range("E15").Value = '  +0.61%  '

# Row 16
$ws.Range("D16").Value = '''3.490.28'
$ws.Range("E16").Value = '  -5.57%  '

# Row 17
$ws.Range("D17").Value = '''61.027.43'
$ws.Range("E17").Value = '  -2.94%  '

# Row 18
$ws.Range("D18").Value = '''3.006.28'
$ws.Range("E18").Value = '  -5.34%  '

# Row 19
$ws.Range("E19").Value = '  -6.76%  '

# Row 20
$ws.Range("D20").Value = '''439.93'
$ws.Range("E20").Value = '  -4.75%  '

# Row 21
$ws.Range("D21").Value = '''13.17'
$ws.Range("E21").Value = '  -5.57%  '

# Row 22
$ws.Range("D22").Value = '''0.664'
$ws.Range("E22").Value = '  -6.80%  '

# Row 23
$ws.Range("D23").Value = '''7.17'
$ws.Range("E23").Value = '  -6.46%  '

# Row 24
$ws.Range("D24").Value = '''12.82'
$ws.Range("E24").Value = '  -4.53%  '

# Row 25
$ws.Range("D25").Value = '''79.27'
$ws.Range("E25").Value = '  -4.97%  '

# Row 26
$ws.Range("D26").Value = '''0.999'
$ws.Range("E26").Value = '  -0.14%  '

# Row 27
$ws.Range("E27").Value = '  -0.02%  '

# Row 28
$ws.Range("D28").Value = '''2.49'
$ws.Range("E28").Value = '  -7.39%  '

# Row 29
$ws.Range("D29").Value = '''7.22'
$ws.Range("E29").Value = '  -7.54%  '

# Row 30
$ws.Range("E30").Value = '  -7.01%  '

# Row 31
$ws.Range("D31").Value = '''25.52'
$ws.Range("E31").Value = '  -6.56%  '

# Row 32
$ws.Range("D32").Value = '''6.11'
$ws.Range("E32").Value = '  -10.05%  '

# Row 33
$ws.Range("D33").Value = '''0.0940'

# Row 34
$ws.Range("E34").Value = '  -4.76%  '

# Row 35
$ws.Range("D35").Value = '''0.958'
$ws.Range("E35").Value = '  -7.93%  '

# Row 36
$ws.Range("D36").Value = '''5.61'
$ws.Range("E36").Value = '  -4.30%  '

# Row 37
$ws.Range("D37").Value = '''50.09'
$ws.Range("E37").Value = '  -1.97%  '

# Row 38
$ws.Range("D38").Value = '0.0₃0672'
$ws.Range("E38").Value = '  -5.65%  '

# Row 39
$ws.Range("E39").Value = '  -6.99%  '

# Row 40
$ws.Range("D40").Value = '''7.79'
$ws.Range("E40").Value = '  -4.03%  '

# Row 41
$ws.Range("D41").Value = '''377.95'
$ws.Range("E41").Value = '  -7.02%  '

# Row 42
$ws.Range("E42").Value = '  -4.35%  '

# Row 43
$ws.Range("D43").Value = '''2.689.78'
$ws.Range("E43").Value = '  -3.76%  '

# Row 44
$ws.Range("E44").Value = '  -8.91%  '

# Row 45
$ws.Range("E45").Value = '  +0.06%  '

# Row 46
$ws.Range("E46").Value = '  -7.10%  '

# Row 47
$ws.Range("B47").Value = 'Monero'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D47").Value = '''120.64'
$ws.Range("E47").Value = '  -2.24%  '

# Row 48
$ws.Range("B48").Value = 'Arweave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D48").Value = '''33.88'
$ws.Range("E48").Value = '  -1.97%  '

# Row 49
$ws.Range("E49").Value = '  -7.84%  '

# Row 50
$ws.Range("E50").Value = '  -4.30%  '

# Row 51
$ws.Range("D51").Value = '''23.37'
$ws.Range("E51").Value = '  -9.22%  '
